$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.567.37"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "1.915.11"
$ws.Range("E3").Value = "  +5.42%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "315.85"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "0.5229"
$ws.Range("E7").Value = "  +4.34%  "
$ws.Range("D8").Value = "0.3967"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").Value = "0.09694"
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("D10").Value = "1.159"
$ws.Range("E10").Value = "  +4.71%  "
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("D12").Value = "6.545"
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").Value = "21.25"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").Value = "1.916.54"
$ws.Range("E14").Value = "  +5.99%  "
$ws.Range("D15").Value = "7.598"
$ws.Range("E15").Value = "  +4.30%  "
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001138"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "94.11"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").Value = "0.06663"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "18.17"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "6.328"
$ws.Range("E22").Value = "  +6.63%  "
$ws.Range("D23").Value = "28.640.31"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("D25").Value = "2.295"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("D26").Value = "2.699"
$ws.Range("E26").Value = "  +11.64%  "
$ws.Range("D27").Value = "2.136.90"
$ws.Range("E27").Value = "  +5.83%  "
$ws.Range("D28").Value = "21.27"
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("D29").Value = "159.83"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("D30").Value = "129.07"
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("D31").Value = "1.107"
$ws.Range("E31").Value = "  +6.78%  "
$ws.Range("D32").Value = "0.1087"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").Value = "5.752"
$ws.Range("E33").Value = "  +3.06%  "
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("E35").Value = "  +10.74%  "
$ws.Range("D36").Value = "0.06784"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("D37").Value = "0.02433"
$ws.Range("E37").Value = "  +4.10%  "
$ws.Range("D38").Value = "1.268"
$ws.Range("E38").Value = "  +6.91%  "
$ws.Range("D39").Value = "0.2226"
$ws.Range("E39").Value = "  +3.87%  "
$ws.Range("D40").Value = "11.88"
$ws.Range("E40").Value = "  +4.98%  "
$ws.Range("D41").Value = "5.115"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6440"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.77%  "
$ws.Range("D43").Value = "1.192"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").Value = "13.57"
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("D46").Value = "0.6091"
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("D48").Value = "1.282"
$ws.Range("D49").Value = "2.037"
$ws.Range("E49").Value = "  +5.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("D51").Value = "1.214"
$ws.Range("E51").Value = "  +2.93%  "
